{"js": "// Edit \"Lista Cerinte\" / \"Lista Specificatii\" document:\n//  1. Shorten the \"Pagin\u0103 cu produse...\" bullet.\n//  2. Shorten the \"Imagini clare \u0219i multiple...\" bullet.\n//  3. Shorten the \"Informa\u021bii despre produs...\" bullet.\n//  4. Shorten the \"Op\u021biuni de plat\u0103 sigure...\" bullet.\n//  5. Remove the \"Posibilitatea de a reseta parola...\" bullet entirely.\n//  6. Remove the \"Responsive Design:\" heading and its single bullet\n//     (\"Website-ul se afi\u0219eaz\u0103 corect \u0219i pe dispozitive mobile.\").\n\nconst body = context.document.body;\n\n// --- Simple text shortenings (search + replace whole paragraph text) ---\nconst replacements = [\n  {\n    find: \"Pagin\u0103 cu produse, cu posibilitate de filtrare \u0219i c\u0103utare.\",\n    replace: \"Pagin\u0103 cu produse.\",\n  },\n  {\n    find: \"Imagini clare \u0219i multiple ale fiec\u0103rui produs.\",\n    replace: \"Imagini multiple ale fiec\u0103rui produs.\",\n  },\n  {\n    find: \"Informa\u021bii despre produs recomand\u0103ri de \u00eengrijire.\",\n    replace: \"Informa\u021bii despre produs.\",\n  },\n  {\n    find: \"Op\u021biuni de plat\u0103 sigure \u0219i variate (card bancar, PayPal, etc.).\",\n    replace: \"Op\u021biuni de plat\u0103 sigure.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const found = body.search(find, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n  found.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n\n// --- Remove whole paragraphs that are no longer part of the requirements list ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst toRemoveTexts = [\n  \"Posibilitatea de a reseta parola \u00een cazul uit\u0103rii.\",\n  \"Responsive Design:\",\n  \"Website-ul se afi\u0219eaz\u0103 corect \u0219i pe dispozitive mobile.\",\n];\n\nconst paragraphsToDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  if (toRemoveTexts.indexOf(text) !== -1) {\n    paragraphsToDelete.push(paragraphs.items[i]);\n  }\n}\n\nif (paragraphsToDelete.length !== toRemoveTexts.length) {\n  throw new Error(\n    \"Expected to find \" + toRemoveTexts.length + \" paragraphs to delete, found \" + paragraphsToDelete.length\n  );\n}\n\nfor (const p of paragraphsToDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# Edit \"Lista Cerinte\" / \"Lista Specificatii\" document:\n#  1. Shorten the \"Pagina cu produse...\" bullet.\n#  2. Shorten the \"Imagini clare si multiple...\" bullet.\n#  3. Shorten the \"Informatii despre produs...\" bullet.\n#  4. Shorten the \"Optiuni de plata sigure...\" bullet.\n#  5. Remove the \"Posibilitatea de a reseta parola...\" bullet entirely.\n#  6. Remove the \"Responsive Design:\" heading and its single bullet\n#     (\"Website-ul se afiseaza corect si pe dispozitive mobile.\").\n\n$d = $word.ActiveDocument\n\n# --- Simple text shortenings (whole-paragraph text swap keeps the\n#     paragraph's existing run formatting/rsid intact, unlike Find.Execute\n#     with a Replacement, which mints a brand-new run) ---\n$replacements = @(\n    @{ Find = \"Pagin\u0103 cu produse, cu posibilitate de filtrare \u0219i c\u0103utare.\"; Replace = \"Pagin\u0103 cu produse.\" },\n    @{ Find = \"Imagini clare \u0219i multiple ale fiec\u0103rui produs.\"; Replace = \"Imagini multiple ale fiec\u0103rui produs.\" },\n    @{ Find = \"Informa\u021bii despre produs recomand\u0103ri de \u00eengrijire.\"; Replace = \"Informa\u021bii despre produs.\" },\n    @{ Find = \"Op\u021biuni de plat\u0103 sigure \u0219i variate (card bancar, PayPal, etc.).\"; Replace = \"Op\u021biuni de plat\u0103 sigure.\" }\n)\n\nforeach ($r in $replacements) {\n    $matched = $false\n    foreach ($p in $d.Paragraphs) {\n        $text = $p.Range.Text.TrimEnd([char]13, [char]10, [char]7)\n        if ($text -eq $r.Find) {\n            $p.Range.Text = $r.Replace\n            $matched = $true\n            break\n        }\n    }\n    if (-not $matched) {\n        throw \"Paragraph not found for: $($r.Find)\"\n    }\n}\n\n# --- Remove whole paragraphs that are no longer part of the requirements list ---\n$toRemoveTexts = @(\n    \"Posibilitatea de a reseta parola \u00een cazul uit\u0103rii.\",\n    \"Responsive Design:\",\n    \"Website-ul se afi\u0219eaz\u0103 corect \u0219i pe dispozitive mobile.\"\n)\n\n# Walk paragraphs back-to-front so deleting one doesn't shift the index of\n# paragraphs we haven't visited yet.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd([char]13, [char]10, [char]7)\n    if ($toRemoveTexts -contains $text) {\n        $p.Range.Delete()\n    }\n}\n"}
